$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.759.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.75%  "

# Row 3
$ws.Range("D3").Value = "'3.019.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.12%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'578.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.89%  "

# Row 6
$ws.Range("D6").Value = "'127.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.98%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").Value = "'3.013.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.26%  "

# Row 9
$ws.Range("D9").Value = "'0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.03%  "

# Row 10
$ws.Range("E10").Value = "  -7.00%  "

# Row 11
$ws.Range("D11").Value = "'5.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.44%  "

# Row 12
$ws.Range("D12").Value = "'0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.45%  "

# Row 13
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.78%  "

# Row 14
$ws.Range("D14").Value = "'32.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.21%  "

# Row 16
$ws.Range("D16").Value = "'3.520.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.92%  "

# Row 17
$ws.Range("D17").Value = "'3.017.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.96%  "

# Row 18
$ws.Range("D18").Value = "'60.667.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.85%  "

# Row 19
$ws.Range("D19").Value = "'6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.45%  "

# Row 20
$ws.Range("D20").Value = "'433.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.13%  "

# Row 21
$ws.Range("D21").Value = "'13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.06%  "

# Row 22
$ws.Range("D22").Value = "'0.668"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.55%  "

# Row 23
$ws.Range("D23").Value = "'7.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.33%  "

# Row 24
$ws.Range("D24").Value = "'12.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.98%  "

# Row 25
$ws.Range("D25").Value = "'79.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.42%  "

# Row 26
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -4.72%  "

# Row 29
$ws.Range("D29").Value = "'7.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.22%  "

# Row 30
$ws.Range("E30").Value = "  -6.96%  "

# Row 31
$ws.Range("D31").Value = "'6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.39%  "

# Row 32
$ws.Range("D32").Value = "'25.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.91%  "

# Row 33
$ws.Range("D33").Value = "'0.0947"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.74%  "

# Row 34
$ws.Range("D34").Value = "'2.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.09%  "

# Row 35
$ws.Range("D35").Value = "'0.963"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.53%  "

# Row 36
$ws.Range("D36").Value = "'5.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.35%  "

# Row 37
$ws.Range("D37").Value = "'50.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.14%  "

# Row 38
$ws.Range("D38").Value = "'0.0₃0676"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.69%  "

# Row 39
$ws.Range("D39").Value = "'8.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "

# Row 40
$ws.Range("E40").Value = "  -7.56%  "

# Row 41
$ws.Range("E41").Value = "  -2.33%  "

# Row 42
$ws.Range("D42").Value = "'387.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.50%  "

# Row 43
$ws.Range("E43").Value = "  -7.87%  "

# Row 44
$ws.Range("D44").Value = "'2.671.98"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "  +0.10%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.238"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.70%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.32%  "

# Row 48
$ws.Range("D48").Value = "'118.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.15%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.01%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'24.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.97%  "

# Row 51
$ws.Range("E51").Value = "  +3.72%  "

